$d = $word.ActiveDocument

# Locate the paragraph that holds the ${FUND1} placeholder.
$target = $null
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -match [regex]::Escape("`${FUND1}")) {
        $target = $p
        $targetIndex = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph containing `${FUND1}"
}

# Insert a brand-new paragraph right after it. InsertParagraphAfter()
# duplicates the source paragraph's formatting (spacing-after 60 twips,
# hanging indent 709/-283, justified), which is exactly what the new
# ${EXTRA} line should use.
$target.Range.InsertParagraphAfter() | Out-Null

# Re-fetch the freshly created paragraph by position and give it its text.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "`${EXTRA}"
